$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-07 Saturday" "2024-12-08 Sunday"

Replace-Text "421÷3=" "816÷9="
Replace-Text "267÷5=" "382÷3="
Replace-Text "643÷9=" "144÷7="
Replace-Text "330÷9=" "562÷6="
Replace-Text "190÷3=" "156÷8="
Replace-Text "324÷3=" "481÷4="
Replace-Text "607÷4=" "403÷5="
Replace-Text "869÷2=" "766÷4="
Replace-Text "648÷2=" "600÷3="
Replace-Text "623÷3=" "560÷9="
Replace-Text "946÷9=" "715÷9="
Replace-Text "845÷4=" "438÷5="
Replace-Text "711÷2=" "935÷7="
Replace-Text "365÷7=" "901÷8="
Replace-Text "437÷3=" "417÷3="
Replace-Text "409÷7=" "613÷4="
Replace-Text "320÷3=" "711÷6="
Replace-Text "359÷7=" "982÷7="
Replace-Text "903÷3=" "673÷2="
Replace-Text "479÷2=" "618÷8="
Replace-Text "446÷2=" "749÷9="
Replace-Text "862÷7=" "438÷2="
Replace-Text "436÷5=" "471÷4="
Replace-Text "432÷2=" "366÷3="
Replace-Text "579÷8=" "674÷7="
